$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# --- Cell values -----------------------------------------------------
# Order matters: it controls the order new entries are appended to the
# shared-strings table, so write in the same sequence the original
# author must have used.

# New "estado" (status) column header on Hoja1
$ws1.Range("B1").Value = "estado"

# Hoja2 becomes the lookup list of possible status values
$ws2.Range("A1").Value = "no comenzado"
$ws2.Range("A2").Value = "en proceso"
$ws2.Range("A3").Value = "terminado"

# Existing two tasks are marked as finished (IIBB percepcion work removed)
$ws1.Range("B2").Value = "terminado"
$ws1.Range("B3").Value = "terminado"

# New task rows
$ws1.Range("A4").Value = "borrar de frontend y backend datos de percepcion"
$ws1.Range("B4").Value = "terminado"
$ws1.Range("A5").Value = "ajustar consultas de cc proveedores listado de fc y de pagos y balance de cuenta"
$ws1.Range("B5").Value = "en proceso"

# --- Column widths -----------------------------------------------------
$ws1.Columns.Item(2).ColumnWidth = 24.0221354166667
$ws2.Columns.Item(1).ColumnWidth = 22.1666666666667

# --- Conditional formatting on Hoja1 column B ---------------------------
$cfRange = $ws1.Range("B1:B1048576")

$fcRed = $cfRange.FormatConditions.Add(1, 3, '"no comenzado"')
$fcRed.Interior.Color = 255

$fcYellow = $cfRange.FormatConditions.Add(1, 3, '"en proceso"')
$fcYellow.Interior.Color = 65535

$fcGreen = $cfRange.FormatConditions.Add(1, 3, '"terminado"')
$fcGreen.Interior.Color = 5296274

# --- Data validation: dropdown list sourced from Hoja2 -------------------
$cfRange.Validation.Add(3, 1, 1, "=Hoja2!`$A`$1:`$A`$3")

# --- Selections ------------------------------------------------------
# Select Hoja2's cursor position first, then Hoja1's last so Hoja1 ends
# up as the active (selected) tab, matching the original workbook.
[void]$ws2.Range("A4").Select()
[void]$ws1.Select()
[void]$ws1.Range("A6").Select()
